# Updated cryptos list on Sun Feb 25 11:30:22 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table,
# and swaps the Monero/Stellar rows (42/43) to their new ranking order.
#
# Note: several Price values are plain decimal-looking strings (e.g.
# "382.90", "0.544"). Setting .Value directly on those makes Excel's COM
# layer auto-coerce them into numbers, which would change the stored cell
# type from text to numeric (not what the source data diff shows - the
# Price column must stay text). To keep them as text without leaving a
# stray NumberFormat behind, we flip the cell to Text format, assign the
# string, then ClearFormats() to drop back to the default (General) style
# while the content stays a literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "51.662.75"
$ws.Range("E2").Value = "  +1.14%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.031.99"
$ws.Range("E3").Value = "  +2.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "382.90"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6 - Solana
Set-TextValue "D6" "102.57"
$ws.Range("E6").Value = "  +0.36%  "

# Row 7 - XRP
Set-TextValue "D7" "0.544"
$ws.Range("E7").Value = "  -0.26%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.38%  "

# Row 10 - Avalanche
Set-TextValue "D10" "36.70"
$ws.Range("E10").Value = "  +0.30%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.01%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +0.96%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.505.61"
$ws.Range("E13").Value = "  +2.37%  "

# Row 14 - Chainlink
Set-TextValue "D14" "18.50"
$ws.Range("E14").Value = "  +0.53%  "

# Row 15 - Polkadot
Set-TextValue "D15" "7.72"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.034.62"
$ws.Range("E16").Value = "  +2.83%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.971"
$ws.Range("E17").Value = "  -3.44%  "

# Row 18 - Uniswap
Set-TextValue "D18" "10.63"
$ws.Range("E18").Value = "  -11.63%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "51.639.04"
$ws.Range("E19").Value = "  +0.98%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -0.46%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "12.44"
$ws.Range("E21").Value = "  +0.52%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  -0.47%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.98"
$ws.Range("E23").Value = "  +0.30%  "

# Row 24 - BitcoinCash (Volume unchanged)
Set-TextValue "D24" "267.36"

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -4.85%  "

# Row 26 - Filecoin
Set-TextValue "D26" "8.42"
$ws.Range("E26").Value = "  +5.87%  "

# Row 27 - RenderToken
Set-TextValue "D27" "7.51"
$ws.Range("E27").Value = "  +6.59%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +3.57%  "

# Row 29 - Dai (unchanged)

# Row 30 - EthereumClassic (Volume unchanged)
Set-TextValue "D30" "26.24"

# Row 31 - Hedera
$ws.Range("E31").Value = "  -1.64%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  -1.93%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  -2.63%  "

# Row 34 - InjectiveProtocol
Set-TextValue "D34" "34.08"
$ws.Range("E34").Value = "  -0.75%  "

# Row 35 - OKB
Set-TextValue "D35" "50.54"
$ws.Range("E35").Value = "  -1.56%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  +2.15%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.11%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -0.07%  "

# Row 39 - TheGraph
Set-TextValue "D39" "0.297"
$ws.Range("E39").Value = "  +7.63%  "

# Row 40 - Celestia
$ws.Range("E40").Value = "  +2.17%  "

# Row 41 - ARBITRUM
Set-TextValue "D41" "1.86"
$ws.Range("E41").Value = "  +1.45%  "

# Row 42 - was Monero, now Stellar
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D42" "0.116"
$ws.Range("E42").Value = "  -0.79%  "

# Row 43 - was Stellar, now Monero
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "127.53"
$ws.Range("E43").Value = "  +2.30%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  +0.92%  "

# Row 45 - NEARProtocol
Set-TextValue "D45" "3.70"
$ws.Range("E45").Value = "  +3.59%  "

# Row 46 - EnergySwap
Set-TextValue "D46" "21.56"
$ws.Range("E46").Value = "  -1.25%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +2.73%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  +2.61%  "

# Row 49 - Maker
Set-TextValue "D49" "2.026.45"
$ws.Range("E49").Value = "  -1.87%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "3.330.36"
$ws.Range("E50").Value = "  +2.48%  "

# Row 51 - WOONetwork
$ws.Range("E51").Value = "  +5.04%  "
